$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

$row = $lo.ListRows.Add()
$r = $row.Range.Row
$ws.Cells.Item($r, 1).Value = "Postman"
$ws.Cells.Item($r, 2).Value = 2021
$ws.Cells.Item($r, 3).Value = 2
$ws.Cells.Item($r, 4).Value = 5
$ws.Cells.Item($r, 5).Value = 8
$ws.Cells.Item($r, 6).Value = 0
$ws.Cells.Item($r, 7).Value = 3
$ws.Cells.Item($r, 8).Formula = "=SUM(C$r`:G$r)"
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 2
$ws.Cells.Item($r, 11).Value = 4

$row = $lo.ListRows.Add()
$r = $row.Range.Row
$ws.Cells.Item($r, 1).Value = "Rippling"
$ws.Cells.Item($r, 2).Value = 2021
$ws.Cells.Item($r, 3).Value = 1
$ws.Cells.Item($r, 4).Value = 4
$ws.Cells.Item($r, 5).Value = 37
$ws.Cells.Item($r, 6).Value = 2
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Formula = "=SUM(C$r`:G$r)"
$ws.Cells.Item($r, 9).Value = 11
$ws.Cells.Item($r, 10).Value = 10
$ws.Cells.Item($r, 11).Value = 1

$row = $lo.ListRows.Add()
$r = $row.Range.Row
$ws.Cells.Item($r, 1).Value = "Salescloud"
$ws.Cells.Item($r, 2).Value = 2021
$ws.Cells.Item($r, 3).Value = 6
$ws.Cells.Item($r, 4).Value = 4
$ws.Cells.Item($r, 5).Value = 40
$ws.Cells.Item($r, 6).Value = 2
$ws.Cells.Item($r, 7).Value = 3
$ws.Cells.Item($r, 8).Formula = "=SUM(C$r`:G$r)"
$ws.Cells.Item($r, 9).Value = 3
$ws.Cells.Item($r, 10).Value = 25
$ws.Cells.Item($r, 11).Value = 4

$ws.Range("K63").Select() | Out-Null
